# carjacking-by-neighborhood-by-month.xlsx — "Add data for 2022-11-18"
#
# The workbook tracks carjackings per Chicago neighborhood (rows) by
# month (columns), where column B holds the running count for the
# current, still-in-progress month. This refresh moves the "through"
# date forward (Nov 09 -> Nov 10) and layers in the incremental counts
# that came in for that day across several historical "November"
# columns as well as the current one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name and the running-month header text both carry the new
# "through" date.
$ws.Name = "Through 2022-11-10"
$ws.Range("B1").Value = "November 2022 (through November 10)"

# Row 4 - New City: add this-month count
$ws.Range("B4").Value = 1

# Row 5 - Garfield Park
$ws.Range("M5").Value = 3     # November 2021: 2 -> 3
$ws.Range("AI5").Value = 3    # November 2019: 2 -> 3
$ws.Range("BE5").Value = 1    # November 2017: new

# Row 9 - Grand Crossing
$ws.Range("M9").Value = 2     # November 2021: 1 -> 2

# Row 11 - Woodlawn
$ws.Range("BE11").Value = 1   # November 2017: new

# Row 16 - Washington Heights
$ws.Range("X16").Value = 3    # November 2020: 1 -> 3
$ws.Range("AI16").Value = 2   # November 2019: 1 -> 2
$ws.Range("AT16").Value = 1   # November 2018: new

# Row 17 - South Shore
$ws.Range("X17").Value = 3    # November 2020: 2 -> 3

# Row 21 - West Town
$ws.Range("BE21").Value = 2   # November 2017: 1 -> 2
$ws.Range("BP21").Value = 1   # November 2016: new

# Row 23 - Albany Park
$ws.Range("X23").Value = 2    # November 2020: 1 -> 2

# Row 26 - Austin
$ws.Range("X26").Value = 4    # November 2020: 3 -> 4
$ws.Range("BE26").Value = 6   # November 2017: 5 -> 6

# Row 31 - Wicker Park
$ws.Range("B31").Value = 1    # this-month: new

# Row 36 - Avondale
$ws.Range("AT36").Value = 2   # November 2018: 1 -> 2

# Row 41 - Roseland
$ws.Range("M41").Value = 3    # November 2021: 2 -> 3

# Row 45 - Calumet Heights
$ws.Range("X45").Value = 1    # November 2020: new

# Row 64 - Gage Park
$ws.Range("AI64").Value = 1   # November 2019: new

# Row 76 - Little Italy, UIC
$ws.Range("BE76").Value = 2   # November 2017: 1 -> 2

# Row 81 - Morgan Park
$ws.Range("AI81").Value = 1   # November 2019: new
